$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Volunteer Details")

# Insert a new blank column at the very start of the sheet; this shifts all
# existing headers/data (old columns A:AJ) one column to the right (new B:AK).
$ws.Columns.Item(1).Insert()

# Populate the newly inserted column A with the new data that was added.
$ws.Range("A1").Value = "No"
$ws.Range("A2").Value = 1

# Move the active selection to match where the author last left the cursor.
$ws.Range("B4").Select() | Out-Null
